# Rename "Sheet1" to "Sprint1" and keep all chart series references
# (which use the sheet name in their formulas) pointing at the correct sheet.

$wb = $excel.ActiveWorkbook

$oldName = "Sheet1"
$newName = "Sprint1"

$ws = $wb.Worksheets.Item($oldName)

# Update any chart series formulas on this sheet so the "Sheet1!" range
# references follow the rename (renaming the worksheet alone does not
# rewrite chart series formulas).
foreach ($chartObj in $ws.ChartObjects()) {
    $chart = $chartObj.Chart
    for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
        $series = $chart.SeriesCollection($i)
        $series.Formula = $series.Formula -replace [regex]::Escape($oldName + "!"), ($newName + "!")
    }
}

# Finally, rename the worksheet itself.
$ws.Name = $newName
